# Update scripts with new TPM values (Fgf6-Fgfr2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) ---
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.09207700000000001
$ws.Range("N2").Value = 0.184154
$ws.Range("O2").Value = 0.0789959771480734
$ws.Range("P2").Value = 0.05545240531440215
$ws.Range("Q2").Value = 0.0004076555713333334
$ws.Range("R2").Value = 0.002445933428
$ws.Range("S2").Value = 0.0789959771480734
$ws.Range("T2").Value = 0.05545240531440215

# --- Row 3 (Target cluster: FAPs) ---
$ws.Range("O3").Value = 0.8491451975864605
$ws.Range("P3").Value = 0.8941052196698643
$ws.Range("Q3").Value = 0.004381979730666666
$ws.Range("R3").Value = 0.039437817576
$ws.Range("S3").Value = 0.8491451975864605
$ws.Range("T3").Value = 0.8941052196698643

# --- Row 4 (Target cluster: MuSCs) ---
$ws.Range("M4").Value = 0.083758
$ws.Range("N4").Value = 0.167516
$ws.Range("O4").Value = 0.07185882526546619
$ws.Range("P4").Value = 0.05044237501573352
$ws.Range("Q4").Value = 0.0003708245853333334
$ws.Range("R4").Value = 0.002224947512
$ws.Range("S4").Value = 0.07185882526546619
$ws.Range("T4").Value = 0.05044237501573352

# --- Row 5 (Target cluster: Neutrophils) removed entirely ---
$ws.Rows.Item(5).Delete()
